# Apply updates to the "想去人数" (F column) values on both the
# "展览" and "全部类型" worksheets, matching the target diff.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F
$updates = @{
    3  = 611
    4  = 2169
    5  = 64
    6  = 12610
    9  = 504
    10 = 464
    12 = 952
    14 = 13987
    19 = 14
    26 = 632
    27 = 5104
    28 = 3
    29 = 254
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
